$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values; regenerated using K instead of Strike# per commit message.
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
